# Updates the cryptos list with refreshed price/volume data.
# Mirrors the commit "Updated cryptos list on Thu Oct 31 23:27:12 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '70.573.06'
$ws.Range("E2").Value2 = '  -2.87%  '
$ws.Range("D3").Value2 = '2.523.51'
$ws.Range("E3").Value2 = '  -5.36%  '
$ws.Range("E4").Value2 = '  -0.03%  '
$c = $ws.Range("D5")
$c.Value2 = "'577.61"
$c.Style = "Normal"
$ws.Range("E5").Value2 = '  -3.39%  '
$c = $ws.Range("D6")
$c.Value2 = "'169.59"
$c.Style = "Normal"
$ws.Range("E6").Value2 = '  -3.51%  '
$ws.Range("E7").Value2 = '  +0.10%  '
$c = $ws.Range("D8")
$c.Value2 = "'0.512"
$c.Style = "Normal"
$ws.Range("E8").Value2 = '  -2.47%  '
$ws.Range("D9").Value2 = '2.522.31'
$ws.Range("E9").Value2 = '  -5.37%  '
$c = $ws.Range("D10")
$c.Value2 = "'0.163"
$c.Style = "Normal"
$ws.Range("E10").Value2 = '  -4.05%  '
$c = $ws.Range("D11")
$c.Value2 = "'0.169"
$c.Style = "Normal"
$ws.Range("E11").Value2 = '  -0.57%  '
$ws.Range("E12").Value2 = '  -3.72%  '
$c = $ws.Range("D13")
$c.Value2 = "'4.84"
$c.Style = "Normal"
$ws.Range("E13").Value2 = '  -3.18%  '
$ws.Range("D14").Value2 = '2.984.24'
$ws.Range("E14").Value2 = '  -5.44%  '
$ws.Range("D15").Value2 = '70.464.50'
$ws.Range("E15").Value2 = '  -2.72%  '
$c = $ws.Range("D16")
$c.Value2 = "'0.0000180"
$c.Style = "Normal"
$ws.Range("E16").Value2 = '  -3.33%  '
$c = $ws.Range("D17")
$c.Value2 = "'25.13"
$c.Style = "Normal"
$ws.Range("E17").Value2 = '  -4.39%  '
$ws.Range("D18").Value2 = '2.516.06'
$ws.Range("E18").Value2 = '  -5.39%  '
$c = $ws.Range("D19")
$c.Value2 = "'11.48"
$c.Style = "Normal"
$ws.Range("E19").Value2 = '  -7.46%  '
$c = $ws.Range("D20")
$c.Value2 = "'7.65"
$c.Style = "Normal"
$ws.Range("E20").Value2 = '  -6.74%  '
$c = $ws.Range("D21")
$c.Value2 = "'360.41"
$c.Style = "Normal"
$ws.Range("E21").Value2 = '  -3.14%  '
$c = $ws.Range("D22")
$c.Value2 = "'3.96"
$c.Style = "Normal"
$ws.Range("E22").Value2 = '  -5.65%  '
$c = $ws.Range("D23")
$c.Value2 = "'1.98"
$c.Style = "Normal"
$ws.Range("E23").Value2 = '  -5.78%  '
$ws.Range("E24").Value2 = '  +0.02%  '
$c = $ws.Range("D25")
$c.Value2 = "'69.45"
$c.Style = "Normal"
$ws.Range("E25").Value2 = '  -3.61%  '
$c = $ws.Range("D26")
$c.Value2 = "'4.07"
$c.Style = "Normal"
$ws.Range("E26").Value2 = '  -6.37%  '
$c = $ws.Range("D27")
$c.Value2 = "'9.15"
$c.Style = "Normal"
$ws.Range("E27").Value2 = '  -7.17%  '
$ws.Range("D28").Value2 = '2.652.51'
$ws.Range("E28").Value2 = '  -4.90%  '
$c = $ws.Range("D29")
$c.Value2 = "'0.985"
$c.Style = "Normal"
$ws.Range("E29").Value2 = '  -1.43%  '
$ws.Range("D30").Value2 = '0.0₃0918'
$ws.Range("E30").Value2 = '  -5.86%  '
$c = $ws.Range("D31")
$c.Value2 = "'7.86"
$c.Style = "Normal"
$ws.Range("E31").Value2 = '  -3.54%  '
$c = $ws.Range("D32")
$c.Value2 = "'485.75"
$c.Style = "Normal"
$ws.Range("E32").Value2 = '  -2.38%  '
$c = $ws.Range("D33")
$c.Value2 = "'1.30"
$c.Style = "Normal"
$ws.Range("E33").Value2 = '  -0.35%  '
$c = $ws.Range("D34")
$c.Value2 = "'1.77"
$c.Style = "Normal"
$ws.Range("E34").Value2 = '  -3.43%  '
$c = $ws.Range("D35")
$c.Value2 = "'0.999"
$c.Style = "Normal"
$ws.Range("E35").Value2 = '  -0.04%  '
$c = $ws.Range("D36")
$c.Value2 = "'156.20"
$c.Style = "Normal"
$ws.Range("E36").Value2 = '  -3.71%  '
$ws.Range("E37").Value2 = '  +1.73%  '
$ws.Range("B38").Value2 = 'EthereumClassic'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D38")
$c.Value2 = "'18.66"
$c.Style = "Normal"
$ws.Range("E38").Value2 = '  -4.53%  '
$ws.Range("B39").Value2 = 'WhiteBITCoin'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D39")
$c.Value2 = "'18.91"
$c.Style = "Normal"
$ws.Range("E39").Value2 = '  -0.21%  '
$c = $ws.Range("D41")
$c.Value2 = "'4.78"
$c.Style = "Normal"
$ws.Range("E41").Value2 = '  -4.71%  '
$c = $ws.Range("D42")
$c.Value2 = "'0.322"
$c.Style = "Normal"
$ws.Range("E42").Value2 = '  -3.54%  '
$c = $ws.Range("D43")
$c.Value2 = "'1.65"
$c.Style = "Normal"
$ws.Range("E43").Value2 = '  -6.89%  '
$c = $ws.Range("D44")
$c.Value2 = "'1.21"
$c.Style = "Normal"
$ws.Range("E44").Value2 = '  -12.85%  '
$c = $ws.Range("D45")
$c.Value2 = "'2.39"
$c.Style = "Normal"
$ws.Range("E45").Value2 = '  -8.03%  '
$c = $ws.Range("D46")
$c.Value2 = "'38.34"
$c.Style = "Normal"
$ws.Range("E46").Value2 = '  -2.48%  '
$c = $ws.Range("D47")
$c.Value2 = "'143.68"
$c.Style = "Normal"
$ws.Range("E47").Value2 = '  -8.59%  '
$c = $ws.Range("D48")
$c.Value2 = "'3.54"
$c.Style = "Normal"
$ws.Range("E48").Value2 = '  -5.46%  '
$c = $ws.Range("D49")
$c.Value2 = "'0.529"
$c.Style = "Normal"
$ws.Range("E49").Value2 = '  -5.48%  '
$ws.Range("E50").Value2 = '  -6.89%  '
$c = $ws.Range("D51")
$c.Value2 = "'0.599"
$c.Style = "Normal"
$ws.Range("E51").Value2 = '  -1.32%  '
